$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column K (the 2020 data column) into a new column L (2021),
# keeping the same formatting as column K.
$ws.Range("K3:K11").Copy() | Out-Null
$ws.Range("L3:L11").PasteSpecial() | Out-Null

# Column L header = 2021
$ws.Range("L4").Value = 2021

# Move the active selection to N2, matching the author's last selection
$ws.Range("N2").Select() | Out-Null
